$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The bold run "run         walk           fly          play              open
# write           swim" becomes two runs: "jump" and the unchanged rest of
# the (tab-look-alike, space padded) verb list, each keeping the same bold
# formatting as before.
$full = $d.Content
$found1 = $full.Find.Execute( `
    "run         walk           fly          play              open         write           swim", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $wordStart = $full.Start
    $rest = "         walk           fly          play              open         write           swim"
    $restLen = $rest.Length

    # Turn "run" into "jump" in place.
    $r1 = $d.Range($wordStart, $wordStart + 3)
    $r1.Text = "jump"

    # Re-apply Bold on the remaining text so it is forced into its own run
    # (same visible formatting, but guarantees a separate <w:r>/<w:t> pair
    # instead of being silently re-merged with the "jump" run).
    $r2 = $d.Range($wordStart + 4, $wordStart + 4 + $restLen)
    $r2.Bold = $false
    $r2.Bold = $true
}

# --- Change 2 -----------------------------------------------------------
# Remove the whole line "I ____________ during breaks at school." (run +
# the line break that introduces it), so the list goes straight from
# "Let's _________________ volleyball." to "Let's ________________ the
# window."
$target = $d.Content
$found2 = $target.Find.Execute("I ____________ during breaks at school.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $delRange = $d.Range($target.Start - 1, $target.End)
    $delRange.Delete()
}
